$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy formatting from the 2021 block (W:Y) into the new 2022 (Z:AB) and 2023 (AC:AE) blocks, row by row ---
for ($row = 4; $row -le 25; $row++) {
    $srcRange = $ws.Range("W" + $row + ":Y" + $row)
    $srcRange.Copy()
    $ws.Range("Z" + $row + ":AB" + $row).PasteSpecial(-4122)
    $srcRange.Copy()
    $ws.Range("AC" + $row + ":AE" + $row).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Step 2: year headers (row 4) ---
$ws.Range("Z4").Value = 2022
$ws.Range("AC4").Value = 2023

# --- Step 3: merge the year header cells ---
$ws.Range("Z4:AB4").Merge()
$ws.Range("AC4:AE4").Merge()

# --- Step 4: sub-header labels (row 5) ---
$ws.Range("Z5").Value = "Both sexes"
$ws.Range("AA5").Value = "Males"
$ws.Range("AB5").Value = "Females"
$ws.Range("AC5").Value = "Both sexes"
$ws.Range("AD5").Value = "Males"
$ws.Range("AE5").Value = "Females"

# --- Step 5: data rows 6-25 ---
$ws.Range("Z6").Value = 499
$ws.Range("AA6").Value = 280
$ws.Range("AB6").Value = 219
$ws.Range("AC6").Value = 453
$ws.Range("AD6").Value = 250
$ws.Range("AE6").Value = 203
$ws.Range("Z7").Value = 4
$ws.Range("AA7").Value = 3
$ws.Range("AB7").Value = 1
$ws.Range("AC7").Value = 2
$ws.Range("AD7").Value = 1
$ws.Range("AE7").Value = 1
$ws.Range("Z8").Value = 1
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 1
$ws.Range("AC8").Value = 1
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 1
$ws.Range("Z9").Value = 0
$ws.Range("AA9").Value = 0
$ws.Range("AB9").Value = 0
$ws.Range("AC9").Value = 1
$ws.Range("AD9").Value = 1
$ws.Range("AE9").Value = 0
$ws.Range("Z10").Value = 0
$ws.Range("AA10").Value = 0
$ws.Range("AB10").Value = 0
$ws.Range("AC10").Value = 0
$ws.Range("AD10").Value = 0
$ws.Range("AE10").Value = 0
$ws.Range("Z11").Value = 0
$ws.Range("AA11").Value = 0
$ws.Range("AB11").Value = 0
$ws.Range("AC11").Value = 0
$ws.Range("AD11").Value = 0
$ws.Range("AE11").Value = 0
$ws.Range("Z12").Value = 2
$ws.Range("AA12").Value = 1
$ws.Range("AB12").Value = 1
$ws.Range("AC12").Value = 1
$ws.Range("AD12").Value = 1
$ws.Range("AE12").Value = 0
$ws.Range("Z13").Value = 2
$ws.Range("AA13").Value = 2
$ws.Range("AB13").Value = 0
$ws.Range("AC13").Value = 3
$ws.Range("AD13").Value = 3
$ws.Range("AE13").Value = 0
$ws.Range("Z14").Value = 4
$ws.Range("AA14").Value = 2
$ws.Range("AB14").Value = 2
$ws.Range("AC14").Value = 2
$ws.Range("AD14").Value = 2
$ws.Range("AE14").Value = 0
$ws.Range("Z15").Value = 6
$ws.Range("AA15").Value = 5
$ws.Range("AB15").Value = 1
$ws.Range("AC15").Value = 7
$ws.Range("AD15").Value = 4
$ws.Range("AE15").Value = 3
$ws.Range("Z16").Value = 6
$ws.Range("AA16").Value = 6
$ws.Range("AB16").Value = 0
$ws.Range("AC16").Value = 8
$ws.Range("AD16").Value = 8
$ws.Range("AE16").Value = 0
$ws.Range("Z17").Value = 12
$ws.Range("AA17").Value = 9
$ws.Range("AB17").Value = 3
$ws.Range("AC17").Value = 11
$ws.Range("AD17").Value = 8
$ws.Range("AE17").Value = 3
$ws.Range("Z18").Value = 16
$ws.Range("AA18").Value = 11
$ws.Range("AB18").Value = 5
$ws.Range("AC18").Value = 15
$ws.Range("AD18").Value = 11
$ws.Range("AE18").Value = 4
$ws.Range("Z19").Value = 30
$ws.Range("AA19").Value = 22
$ws.Range("AB19").Value = 8
$ws.Range("AC19").Value = 25
$ws.Range("AD19").Value = 20
$ws.Range("AE19").Value = 5
$ws.Range("Z20").Value = 43
$ws.Range("AA20").Value = 31
$ws.Range("AB20").Value = 12
$ws.Range("AC20").Value = 33
$ws.Range("AD20").Value = 22
$ws.Range("AE20").Value = 11
$ws.Range("Z21").Value = 45
$ws.Range("AA21").Value = 31
$ws.Range("AB21").Value = 14
$ws.Range("AC21").Value = 51
$ws.Range("AD21").Value = 36
$ws.Range("AE21").Value = 15
$ws.Range("Z22").Value = 67
$ws.Range("AA22").Value = 37
$ws.Range("AB22").Value = 30
$ws.Range("AC22").Value = 72
$ws.Range("AD22").Value = 38
$ws.Range("AE22").Value = 34
$ws.Range("Z23").Value = 45
$ws.Range("AA23").Value = 22
$ws.Range("AB23").Value = 23
$ws.Range("AC23").Value = 54
$ws.Range("AD23").Value = 23
$ws.Range("AE23").Value = 31
$ws.Range("Z24").Value = 120
$ws.Range("AA24").Value = 51
$ws.Range("AB24").Value = 69
$ws.Range("AC24").Value = 71
$ws.Range("AD24").Value = 34
$ws.Range("AE24").Value = 37
$ws.Range("Z25").Value = 96
$ws.Range("AA25").Value = 47
$ws.Range("AB25").Value = 49
$ws.Range("AC25").Value = 96
$ws.Range("AD25").Value = 38
$ws.Range("AE25").Value = 58

# --- Step 6: selection state to match target sheetView ---
$ws.Range("A4:A5").Select()
